{"js": "// Office.js (Word JavaScript API) script.\n// Replaces each arithmetic-problem cell text with its new value, per the\n// commit diff (old expression -> new expression). Every old string occurs\n// exactly once in the document body, so a simple search+replace per pair\n// is safe and order-independent.\n\nconst replacements = [\n  [\"94-59=35\", \"1+31=32\"],\n  [\"6+20=26\", \"81-45=36\"],\n  [\"83-79=4\", \"35+58=93\"],\n  [\"60-13=47\", \"84-38=46\"],\n  [\"14+61=75\", \"28+18=46\"],\n  [\"32+43=75\", \"83-20=63\"],\n  [\"41+32=73\", \"62+35=97\"],\n  [\"31+38=69\", \"11+46=57\"],\n  [\"53+29=82\", \"50-44=6\"],\n  [\"26+57=83\", \"50-46=4\"],\n  [\"4+79=83\", \"20+31=51\"],\n  [\"96-61=35\", \"14+4=18\"],\n  [\"44-16=28\", \"3+71=74\"],\n  [\"26+54=80\", \"40+0=40\"],\n  [\"36+51=87\", \"7+91=98\"],\n  [\"10+78=88\", \"81-43=38\"],\n  [\"62+29=91\", \"71-47=24\"],\n  [\"36+46=82\", \"71+0=71\"],\n  [\"81+3=84\", \"36-8=28\"],\n  [\"97-91=6\", \"75-57=18\"],\n  [\"24+61=85\", \"66-2=64\"],\n  [\"11+24=35\", \"78-1=77\"],\n  [\"41+48=89\", \"35+43=78\"],\n  [\"12+12=24\", \"5-4=1\"],\n  [\"28-13=15\", \"28-9=19\"],\n  [\"57-6=51\", \"59-57=2\"],\n  [\"12+48=60\", \"63+27=90\"],\n  [\"11+33=44\", \"66-17=49\"],\n  [\"23+24=47\", \"87-40=47\"],\n  [\"73-70=3\", \"22+6=28\"],\n  [\"53+46=99\", \"76-17=59\"],\n  [\"97-75=22\", \"21-3=18\"],\n  [\"82-20=62\", \"59+6=65\"],\n  [\"91-51=40\", \"97-34=63\"],\n  [\"48+41=89\", \"94+4=98\"],\n  [\"38-27=11\", \"4+65=69\"],\n  [\"23-14=9\", \"65-18=47\"],\n  [\"84-46=38\", \"56-21=35\"],\n  [\"17+38=55\", \"35+5=40\"],\n  [\"97-6=91\", \"45+25=70\"],\n  [\"77-25=52\", \"38+26=64\"],\n  [\"74-72=2\", \"59-18=41\"],\n  [\"42-17=25\", \"82+10=92\"],\n  [\"60-31=29\", \"50-33=17\"],\n  [\"55-26=29\", \"13+54=67\"],\n  [\"29-3=26\", \"22+55=77\"],\n  [\"9+76=85\", \"26+29=55\"],\n  [\"58+39=97\", \"55-15=40\"],\n  [\"6+43=49\", \"12-1=11\"],\n  [\"26+39=65\", \"69+28=97\"],\n  [\"26+51=77\", \"49-31=18\"],\n  [\"93-35=58\", \"62-51=11\"],\n  [\"96-57=39\", \"33-32=1\"],\n  [\"85-5=80\", \"42+31=73\"],\n  [\"23+39=62\", \"85-28=57\"],\n  [\"19+15=34\", \"6+49=55\"],\n  [\"62-39=23\", \"92-53=39\"],\n  [\"95-65=30\", \"85-62=23\"],\n  [\"44-9=35\", \"90-77=13\"],\n  [\"45-14=31\", \"68-63=5\"],\n  [\"76-56=20\", \"77-20=57\"],\n  [\"47-7=40\", \"45+13=58\"],\n  [\"61-51=10\", \"39+11=50\"],\n  [\"8+32=40\", \"74-46=28\"],\n  [\"46-39=7\", \"5+73=78\"],\n  [\"31+66=97\", \"65+21=86\"],\n  [\"66-46=20\", \"28-6=22\"],\n  [\"80+7=87\", \"49-21=28\"],\n  [\"14+14=28\", \"1+16=17\"],\n  [\"99-55=44\", \"27+11=38\"],\n  [\"42-13=29\", \"27+15=42\"],\n  [\"90-16=74\", \"17+60=77\"],\n  [\"35-25=10\", \"48-44=4\"],\n  [\"49+16=65\", \"57-30=27\"],\n  [\"97-81=16\", \"65-24=41\"],\n  [\"52+41=93\", \"47-24=23\"],\n  [\"12+31=43\", \"0+90=90\"],\n  [\"52-33=19\", \"32+41=73\"],\n  [\"66-36=30\", \"35+3=38\"],\n  [\"38+38=76\", \"81-5=76\"],\n  [\"64-52=12\", \"29+50=79\"],\n  [\"51-32=19\", \"3+37=40\"],\n  [\"85-24=61\", \"87-41=46\"],\n  [\"99-67=32\", \"22+11=33\"],\n  [\"3+64=67\", \"62-27=35\"],\n  [\"89-32=57\", \"69-54=15\"],\n  [\"22-17=5\", \"60+20=80\"],\n  [\"9+74=83\", \"73-32=41\"],\n  [\"60-58=2\", \"75+0=75\"],\n  [\"97-65=32\", \"56+21=77\"],\n  [\"6+90=96\", \"81-45=36\"],\n  [\"34-9=25\", \"32+63=95\"],\n  [\"28-27=1\", \"18+55=73\"],\n  [\"80+1=81\", \"88-78=10\"],\n  [\"86-22=64\", \"94-15=79\"],\n  [\"99-61=38\", \"78-26=52\"],\n  [\"69-49=20\", \"18+14=32\"],\n  [\"44-29=15\", \"1+70=71\"],\n  [\"42-6=36\", \"12+55=67\"],\n  [\"90-39=51\", \"8+31=39\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Replaces each arithmetic-problem cell text with its new value, per the\n# commit diff (old expression -> new expression). Each old string occurs\n# exactly once in the document body.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"94-59=35\", \"1+31=32\"),\n    @(\"6+20=26\", \"81-45=36\"),\n    @(\"83-79=4\", \"35+58=93\"),\n    @(\"60-13=47\", \"84-38=46\"),\n    @(\"14+61=75\", \"28+18=46\"),\n    @(\"32+43=75\", \"83-20=63\"),\n    @(\"41+32=73\", \"62+35=97\"),\n    @(\"31+38=69\", \"11+46=57\"),\n    @(\"53+29=82\", \"50-44=6\"),\n    @(\"26+57=83\", \"50-46=4\"),\n    @(\"4+79=83\", \"20+31=51\"),\n    @(\"96-61=35\", \"14+4=18\"),\n    @(\"44-16=28\", \"3+71=74\"),\n    @(\"26+54=80\", \"40+0=40\"),\n    @(\"36+51=87\", \"7+91=98\"),\n    @(\"10+78=88\", \"81-43=38\"),\n    @(\"62+29=91\", \"71-47=24\"),\n    @(\"36+46=82\", \"71+0=71\"),\n    @(\"81+3=84\", \"36-8=28\"),\n    @(\"97-91=6\", \"75-57=18\"),\n    @(\"24+61=85\", \"66-2=64\"),\n    @(\"11+24=35\", \"78-1=77\"),\n    @(\"41+48=89\", \"35+43=78\"),\n    @(\"12+12=24\", \"5-4=1\"),\n    @(\"28-13=15\", \"28-9=19\"),\n    @(\"57-6=51\", \"59-57=2\"),\n    @(\"12+48=60\", \"63+27=90\"),\n    @(\"11+33=44\", \"66-17=49\"),\n    @(\"23+24=47\", \"87-40=47\"),\n    @(\"73-70=3\", \"22+6=28\"),\n    @(\"53+46=99\", \"76-17=59\"),\n    @(\"97-75=22\", \"21-3=18\"),\n    @(\"82-20=62\", \"59+6=65\"),\n    @(\"91-51=40\", \"97-34=63\"),\n    @(\"48+41=89\", \"94+4=98\"),\n    @(\"38-27=11\", \"4+65=69\"),\n    @(\"23-14=9\", \"65-18=47\"),\n    @(\"84-46=38\", \"56-21=35\"),\n    @(\"17+38=55\", \"35+5=40\"),\n    @(\"97-6=91\", \"45+25=70\"),\n    @(\"77-25=52\", \"38+26=64\"),\n    @(\"74-72=2\", \"59-18=41\"),\n    @(\"42-17=25\", \"82+10=92\"),\n    @(\"60-31=29\", \"50-33=17\"),\n    @(\"55-26=29\", \"13+54=67\"),\n    @(\"29-3=26\", \"22+55=77\"),\n    @(\"9+76=85\", \"26+29=55\"),\n    @(\"58+39=97\", \"55-15=40\"),\n    @(\"6+43=49\", \"12-1=11\"),\n    @(\"26+39=65\", \"69+28=97\"),\n    @(\"26+51=77\", \"49-31=18\"),\n    @(\"93-35=58\", \"62-51=11\"),\n    @(\"96-57=39\", \"33-32=1\"),\n    @(\"85-5=80\", \"42+31=73\"),\n    @(\"23+39=62\", \"85-28=57\"),\n    @(\"19+15=34\", \"6+49=55\"),\n    @(\"62-39=23\", \"92-53=39\"),\n    @(\"95-65=30\", \"85-62=23\"),\n    @(\"44-9=35\", \"90-77=13\"),\n    @(\"45-14=31\", \"68-63=5\"),\n    @(\"76-56=20\", \"77-20=57\"),\n    @(\"47-7=40\", \"45+13=58\"),\n    @(\"61-51=10\", \"39+11=50\"),\n    @(\"8+32=40\", \"74-46=28\"),\n    @(\"46-39=7\", \"5+73=78\"),\n    @(\"31+66=97\", \"65+21=86\"),\n    @(\"66-46=20\", \"28-6=22\"),\n    @(\"80+7=87\", \"49-21=28\"),\n    @(\"14+14=28\", \"1+16=17\"),\n    @(\"99-55=44\", \"27+11=38\"),\n    @(\"42-13=29\", \"27+15=42\"),\n    @(\"90-16=74\", \"17+60=77\"),\n    @(\"35-25=10\", \"48-44=4\"),\n    @(\"49+16=65\", \"57-30=27\"),\n    @(\"97-81=16\", \"65-24=41\"),\n    @(\"52+41=93\", \"47-24=23\"),\n    @(\"12+31=43\", \"0+90=90\"),\n    @(\"52-33=19\", \"32+41=73\"),\n    @(\"66-36=30\", \"35+3=38\"),\n    @(\"38+38=76\", \"81-5=76\"),\n    @(\"64-52=12\", \"29+50=79\"),\n    @(\"51-32=19\", \"3+37=40\"),\n    @(\"85-24=61\", \"87-41=46\"),\n    @(\"99-67=32\", \"22+11=33\"),\n    @(\"3+64=67\", \"62-27=35\"),\n    @(\"89-32=57\", \"69-54=15\"),\n    @(\"22-17=5\", \"60+20=80\"),\n    @(\"9+74=83\", \"73-32=41\"),\n    @(\"60-58=2\", \"75+0=75\"),\n    @(\"97-65=32\", \"56+21=77\"),\n    @(\"6+90=96\", \"81-45=36\"),\n    @(\"34-9=25\", \"32+63=95\"),\n    @(\"28-27=1\", \"18+55=73\"),\n    @(\"80+1=81\", \"88-78=10\"),\n    @(\"86-22=64\", \"94-15=79\"),\n    @(\"99-61=38\", \"78-26=52\"),\n    @(\"69-49=20\", \"18+14=32\"),\n    @(\"44-29=15\", \"1+70=71\"),\n    @(\"42-6=36\", \"12+55=67\"),\n    @(\"90-39=51\", \"8+31=39\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n\n"}
